$d = $word.ActiveDocument

$pairs = @(
    @("772÷7=", "882÷5="),
    @("775÷2=", "194÷8="),
    @("915÷4=", "882÷3="),
    @("737÷9=", "434÷9="),
    @("761÷4=", "524÷2="),
    @("963÷2=", "960÷5="),
    @("431÷8=", "766÷3="),
    @("183÷9=", "871÷8="),
    @("503÷8=", "728÷8="),
    @("623÷5=", "250÷4="),
    @("925÷8=", "424÷6="),
    @("680÷5=", "363÷9="),
    @("131÷9=", "937÷8="),
    @("934÷5=", "260÷7="),
    @("499÷5=", "723÷6="),
    @("348÷6=", "871÷4="),
    @("717÷3=", "564÷6="),
    @("482÷3=", "269÷5="),
    @("584÷5=", "753÷7="),
    @("141÷2=", "133÷4="),
    @("599÷2=", "699÷5="),
    @("389÷3=", "915÷6="),
    @("987÷9=", "339÷8="),
    @("729÷2=", "203÷2="),
    @("187÷8=", "978÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
